$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.127733
$ws.Range("H2").Value = 78.383199
$ws.Range("I2").Value = 0.2666992864894373
$ws.Range("J2").Value = 0.2666992864894374
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.03508633333333333
$ws.Range("N2").Value = 0.105259
$ws.Range("O2").Value = 0.004489537393262644
$ws.Range("P2").Value = 0.004489537393262644
$ws.Range("Q2").Value = 0.9167263492823333
$ws.Range("R2").Value = 8.250537143540999
$ws.Range("S2").Value = 0.001197356419450795
$ws.Range("T2").Value = 0.001197356419450796
$ws.Range("G3").Value = 26.127733
$ws.Range("H3").Value = 78.383199
$ws.Range("I3").Value = 0.2666992864894373
$ws.Range("J3").Value = 0.2666992864894374
$ws.Range("M3").Value = 4.911922333333334
$ws.Range("O3").Value = 0.6285142074777995
$ws.Range("P3").Value = 0.6285142074777995
$ws.Range("Q3").Value = 128.3373952420704
$ws.Range("R3").Value = 1155.036557178633
$ws.Range("S3").Value = 0.1676242906828033
$ws.Range("T3").Value = 0.1676242906828033
$ws.Range("G4").Value = 26.127733
$ws.Range("H4").Value = 78.383199
$ws.Range("I4").Value = 0.2666992864894373
$ws.Range("J4").Value = 0.2666992864894374
$ws.Range("M4").Value = 2.868124666666667
$ws.Range("N4").Value = 8.604374
$ws.Range("O4").Value = 0.3669962551289379
$ws.Range("P4").Value = 0.3669962551289379
$ws.Range("Q4").Value = 74.93759550138067
$ws.Range("R4").Value = 674.438359512426
$ws.Range("S4").Value = 0.09787763938718322
$ws.Range("T4").Value = 0.09787763938718325
$ws.Range("I5").Value = 0.2440410104700376
$ws.Range("J5").Value = 0.2440410104700377
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.03508633333333333
$ws.Range("N5").Value = 0.105259
$ws.Range("O5").Value = 0.004489537393262644
$ws.Range("P5").Value = 0.004489537393262644
$ws.Range("Q5").Value = 0.8388429813524443
$ws.Range("R5").Value = 7.549586832171999
$ws.Range("S5").Value = 0.001095631241994834
$ws.Range("T5").Value = 0.001095631241994834
$ws.Range("I6").Value = 0.2440410104700376
$ws.Range("J6").Value = 0.2440410104700377
$ws.Range("M6").Value = 4.911922333333334
$ws.Range("O6").Value = 0.6285142074777995
$ws.Range("P6").Value = 0.6285142074777995
$ws.Range("Q6").Value = 117.4340885130485
$ws.Range("S6").Value = 0.153383242287657
$ws.Range("T6").Value = 0.1533832422876571
$ws.Range("I7").Value = 0.2440410104700376
$ws.Range("J7").Value = 0.2440410104700377
$ws.Range("M7").Value = 2.868124666666667
$ws.Range("N7").Value = 8.604374
$ws.Range("O7").Value = 0.3669962551289379
$ws.Range("P7").Value = 0.3669962551289379
$ws.Range("Q7").Value = 68.57103657484355
$ws.Range("R7").Value = 617.1393291735919
$ws.Range("S7").Value = 0.08956213694038571
$ws.Range("T7").Value = 0.08956213694038573
$ws.Range("G8").Value = 47.93131266666666
$ws.Range("H8").Value = 143.793938
$ws.Range("I8").Value = 0.489259703040525
$ws.Range("J8").Value = 0.4892597030405251
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.03508633333333333
$ws.Range("N8").Value = 0.105259
$ws.Range("O8").Value = 0.004489537393262644
$ws.Range("P8").Value = 0.004489537393262644
$ws.Range("Q8").Value = 1.681734013326889
$ws.Range("R8").Value = 15.135606119942
$ws.Range("S8").Value = 0.002196549731817014
$ws.Range("T8").Value = 0.002196549731817014
$ws.Range("G9").Value = 47.93131266666666
$ws.Range("H9").Value = 143.793938
$ws.Range("I9").Value = 0.489259703040525
$ws.Range("J9").Value = 0.4892597030405251
$ws.Range("M9").Value = 4.911922333333334
$ws.Range("O9").Value = 0.6285142074777995
$ws.Range("P9").Value = 0.6285142074777995
$ws.Range("Q9").Value = 235.4348851533829
$ws.Range("R9").Value = 2118.913966380446
$ws.Range("S9").Value = 0.3075066745073391
$ws.Range("T9").Value = 0.3075066745073392
$ws.Range("G10").Value = 47.93131266666666
$ws.Range("H10").Value = 143.793938
$ws.Range("I10").Value = 0.489259703040525
$ws.Range("J10").Value = 0.4892597030405251
$ws.Range("M10").Value = 2.868124666666667
$ws.Range("N10").Value = 8.604374
$ws.Range("O10").Value = 0.3669962551289379
$ws.Range("P10").Value = 0.3669962551289379
$ws.Range("Q10").Value = 137.4729801649791
$ws.Range("R10").Value = 1237.256821484812
$ws.Range("S10").Value = 0.1795564788013689
$ws.Range("T10").Value = 0.1795564788013689
